$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.523.23'
$ws.Range('E2').Value = '  +3.15%  '
$ws.Range('D3').Value = '1.843.54'
$ws.Range('E3').Value = '  +2.43%  '
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '232.62'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('E5').Value = '  +3.66%  '
$ws.Range('E6').Value = '  +2.78%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '44.42'
$ws.Range('D8').NumberFormat = "General"
$ws.Range('E8').Value = '  +13.78%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.311'
$ws.Range('D9').NumberFormat = "General"
$ws.Range('E9').Value = '  +8.22%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0698'
$ws.Range('D10').NumberFormat = "General"
$ws.Range('E10').Value = '  +4.49%  '
$ws.Range('E11').Value = '  +2.60%  '
$ws.Range('D12').Value = '2.111.11'
$ws.Range('E12').Value = '  +2.51%  '
$ws.Range('D13').Value = '1.844.03'
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.31'
$ws.Range('D14').NumberFormat = "General"
$ws.Range('E14').Value = '  +4.06%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.675'
$ws.Range('D15').NumberFormat = "General"
$ws.Range('E15').Value = '  +7.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.72'
$ws.Range('D16').NumberFormat = "General"
$ws.Range('E16').Value = '  +8.42%  '
$ws.Range('D17').Value = '35.517.91'
$ws.Range('E17').Value = '  +3.20%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '70.48'
$ws.Range('D18').NumberFormat = "General"
$ws.Range('E18').Value = '  +3.74%  '
$ws.Range('E19').Value = '  +5.15%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '244.05'
$ws.Range('D20').NumberFormat = "General"
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.10'
$ws.Range('D21').NumberFormat = "General"
$ws.Range('E21').Value = '  +9.44%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.62'
$ws.Range('D22').NumberFormat = "General"
$ws.Range('E22').Value = '  +13.58%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.25'
$ws.Range('D24').NumberFormat = "General"
$ws.Range('E24').Value = '  +4.19%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '171.44'
$ws.Range('D25').NumberFormat = "General"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.00'
$ws.Range('D26').NumberFormat = "General"
$ws.Range('E26').Value = '  +4.52%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.84'
$ws.Range('D27').NumberFormat = "General"
$ws.Range('E28').Value = '  +1.47%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.56'
$ws.Range('D29').NumberFormat = "General"
$ws.Range('E29').Value = '  +28.34%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').Value = '3.328.47'
$ws.Range('E31').Value = '  +36.99%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0555'
$ws.Range('D32').NumberFormat = "General"
$ws.Range('E32').Value = '  +8.40%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.12'
$ws.Range('D33').NumberFormat = "General"
$ws.Range('E33').Value = '  +7.89%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.94'
$ws.Range('D34').NumberFormat = "General"
$ws.Range('E34').Value = '  +5.54%  '
$ws.Range('E35').Value = '  +2.28%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '95.85'
$ws.Range('D36').NumberFormat = "General"
$ws.Range('E36').Value = '  +17.40%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.694'
$ws.Range('D37').NumberFormat = "General"
$ws.Range('E37').Value = '  +8.70%  '
$ws.Range('E38').Value = '  +8.82%  '
$ws.Range('D39').Value = '1.346.83'
$ws.Range('E39').Value = '  +3.11%  '
$ws.Range('E40').Value = '  +5.88%  '
$ws.Range('E41').Value = '  +6.37%  '
$ws.Range('E42').Value = '  +7.85%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '15.23'
$ws.Range('D43').NumberFormat = "General"
$ws.Range('E43').Value = '  +8.52%  '
$ws.Range('E44').Value = '  +3.22%  '
$ws.Range('E45').Value = '  +0.81%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.81'
$ws.Range('D46').NumberFormat = "General"
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.27'
$ws.Range('D47').NumberFormat = "General"
$ws.Range('E47').Value = '  +9.86%  '
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('D49').Value = '2.017.18'
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '102.49'
$ws.Range('D51').NumberFormat = "General"
$ws.Range('E51').Value = '  +0.72%  '
